$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "cancer_name" (column F) and "icd-10" (column G) columns were exported
# with fixed-width trailing padding baked into the text (e.g. "Lip" followed
# by many spaces). This strips that trailing whitespace, equivalent to
# Python's str.strip(), from every data row in those two columns.

$lastRow = $ws.Cells($ws.Rows.Count, 6).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cF = $ws.Cells.Item($r, 6)
    $fVal = $cF.Value2
    if ($fVal -ne $null) {
        $trimmedF = $fVal.Trim()
        if ($trimmedF -ne $fVal) {
            $cF.Value = $trimmedF
        }
    }

    $cG = $ws.Cells.Item($r, 7)
    $gVal = $cG.Value2
    if ($gVal -ne $null) {
        $trimmedG = $gVal.Trim()
        if ($trimmedG -ne $gVal) {
            $cG.Value = $trimmedG
        }
    }
}
